# Generate Report for Handback
# This script updates the localization-status report after a handback:
#  - the zh-cn / de-de "Status" cells move from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - the "Latest Target File" / "Latest Handback File" / "Latest Handback
#    DateTime" columns on the zh-cn and de-de sheets are filled in
#  - column widths are refreshed to fit the new, longer content

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E and F) for both rows
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn / de-de sheets: Status column (C) for both rows
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn / de-de sheets: Latest Target File (I) now links to the source
# markdown file, same as column A already does.
# ---------------------------------------------------------------------
$mdFileName     = "f8f155ac-b6d9-4825-800a-3cac39ded6bb.md"
$mdHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3a99899f34e516c2614f6d18f9ead28b6621921/e2e/f8f155ac-b6d9-4825-800a-3cac39ded6bb.md"
$hyperlinkColor = 15570276   # long-form BGR of RGB(0x64,0x95,0xED) - matches the existing hyperlink style

function Add-TargetFileLink {
    param($ws, [string]$cellAddr)

    $ws.Hyperlinks.Add($ws.Range($cellAddr), $mdHyperlinkUrl, $null, $mdFileName, $mdFileName)
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = $hyperlinkColor
}

Add-TargetFileLink $wsZhCn "I2"
Add-TargetFileLink $wsZhCn "I3"
Add-TargetFileLink $wsDeDe "I2"
Add-TargetFileLink $wsDeDe "I3"

# ---------------------------------------------------------------------
# zh-cn / de-de sheets: Latest Handback File (J) - the generated xliff
# ---------------------------------------------------------------------
$wsZhCn.Range("J2").Value = "f8f155ac-b6d9-4825-800a-3cac39ded6bb.8fe65fe2b16217659ca8b752011a1fbee6de9245.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "f8f155ac-b6d9-4825-800a-3cac39ded6bb.8fe65fe2b16217659ca8b752011a1fbee6de9245.zh-cn.xlf"
$wsDeDe.Range("J2").Value = "f8f155ac-b6d9-4825-800a-3cac39ded6bb.8fe65fe2b16217659ca8b752011a1fbee6de9245.de-de.xlf"
$wsDeDe.Range("J3").Value = "f8f155ac-b6d9-4825-800a-3cac39ded6bb.8fe65fe2b16217659ca8b752011a1fbee6de9245.de-de.xlf"

# ---------------------------------------------------------------------
# zh-cn / de-de sheets: Latest Handback DateTime (K)
# ---------------------------------------------------------------------
$wsZhCn.Range("K2").Value = "2016-08-16 23:01:41"
$wsZhCn.Range("K3").Value = "2016-08-16 23:01:41"
$wsDeDe.Range("K2").Value = "2016-08-16 23:01:48"
$wsDeDe.Range("K3").Value = "2016-08-16 23:01:48"

# ---------------------------------------------------------------------
# Refresh column widths that depend on the new, longer cell content.
# ColumnWidth is padded by ~0.8333 (5/6) characters when stored, so the
# requested width is the desired stored width minus that padding.
# ---------------------------------------------------------------------
function Set-StoredColumnWidth {
    param($column, [double]$storedWidth)
    $column.ColumnWidth = $storedWidth - (5.0 / 6.0)
}

Set-StoredColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527   # E - zh-cn status
Set-StoredColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527   # F - de-de status

Set-StoredColumnWidth $wsZhCn.Columns.Item(3) 29.9777047293527       # C - Status
Set-StoredColumnWidth $wsZhCn.Columns.Item(9) 40                     # I - Latest Target File
Set-StoredColumnWidth $wsZhCn.Columns.Item(10) 40                    # J - Latest Handback File

Set-StoredColumnWidth $wsDeDe.Columns.Item(3) 29.9777047293527       # C - Status
Set-StoredColumnWidth $wsDeDe.Columns.Item(9) 40                     # I - Latest Target File
Set-StoredColumnWidth $wsDeDe.Columns.Item(10) 40                    # J - Latest Handback File
